$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 701, shifting all subsequent rows down by one.
$ws.Rows("701:701").Insert()

# Fill in the new row's data. Column A holds date-like text, so use a
# leading apostrophe to force text entry (avoiding auto date conversion),
# then restore the cell's style to match the surrounding plain cells.
$ws.Range("A701").Value = "'2026/01/24"
$ws.Range("A701").Style = "Normal"
$ws.Range("B701").Value = "土"
$ws.Range("C701").Value = 3
$ws.Range("D701").Value = 201
